$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: merge the paragraph ending in "...(Quadro 2) " with the
# following paragraph that only contains "[IMG] . " - the result is a
# single paragraph whose sole run holds the concatenated text, and the
# now-redundant paragraph (and its mark) disappear.
# ------------------------------------------------------------------
$quadroIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("Quadro 2)")) {
        $quadroIndex = $i
        break
    }
}

if ($quadroIndex -eq -1) {
    throw "Could not find the paragraph containing 'Quadro 2)'"
}

$quadroPara = $d.Paragraphs.Item($quadroIndex)
$imgPara = $quadroPara.Next()

# Sanity check: the following paragraph should start with "[IMG] . "
# (NB: use .StartsWith, not -like, since [ and ] are -like wildcard
# metacharacters and would otherwise be parsed as a character class.)
if (-not ($imgPara.Range.Text.StartsWith("[IMG] . "))) {
    throw "Unexpected paragraph following the Quadro 2 paragraph: [$($imgPara.Range.Text)]"
}

# Text of that paragraph, without its trailing paragraph mark.
$imgText = $d.Range($imgPara.Range.Start, $imgPara.Range.End - 1).Text

# Range covering the Quadro-2 paragraph's text, excluding its paragraph
# mark, so the appended text lands inside the existing run.
$quadroTextRange = $d.Range($quadroPara.Range.Start, $quadroPara.Range.End - 1)
$quadroTextRange.InsertAfter($imgText)

# Re-fetch the (now stand-alone, to-be-removed) "[IMG] . " paragraph and
# delete it fully, paragraph mark included - this merges it away, leaving
# a single paragraph with one run holding the combined text.
$quadroPara2 = $d.Paragraphs.Item($quadroIndex)
$imgPara2 = $quadroPara2.Next()
$imgPara2.Range.Delete()

# ------------------------------------------------------------------
# Edit 2: the final "w$ D" answer-key paragraph (last occurrence in the
# document) gains a trailing space on its existing run, an empty bold
# run, and a new bold run containing "a".
# ------------------------------------------------------------------
$wdIndex = -1
$count2 = $d.Paragraphs.Count
for ($i = 1; $i -le $count2; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $d.Range($p.Range.Start, $p.Range.End - 1).Text
    if ($t -eq "w`$ D") {
        $wdIndex = $i
    }
}

if ($wdIndex -eq -1) {
    throw "Could not find a paragraph containing exactly 'w`$ D'"
}

$wdPara = $d.Paragraphs.Item($wdIndex)
$wdRange = $d.Range($wdPara.Range.Start, $wdPara.Range.End - 1)

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">w$ D </w:t></w:r><w:r><w:rPr><w:b/></w:rPr></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>a</w:t></w:r></w:p>'
$wdRange.InsertXML($newXml)
